$d = $word.ActiveDocument

# The document currently ends with a paragraph:
#   "5. To guarantee 3 pairs you need to pick out all 20 socks." + (hidden) _GoBack bookmark
# We need to:
#   1. Close that paragraph right after the text (removing the trailing _GoBack bookmark from it).
#   2. Add a brand-new paragraph after it containing:
#        "4. A. Each solution somewhat meets the goals the only that works for both questions is
#         that you need " (first run, carries a lastRenderedPageBreak marker)
#        the _GoBack bookmark (empty, start immediately followed by end)
#        "to pick 20 socks to guarantee 1 or 3 pairs." (second run)

# The existing _GoBack bookmark sits at the very end of the document's content. Adding/placing a
# bookmark exactly at the absolute last character position of the whole story is unreliable, so
# remove the old one first; it gets recreated (in its new home) later in this script.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# Paragraph that currently ends the document: "...pick out all 20 socks."
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)
$lastRange = $lastPara.Range

# Collapse to just after the final period, i.e. right before the paragraph mark, and split the
# paragraph there (equivalent to placing the cursor there and pressing Enter).
$splitPoint = $d.Range($lastRange.End - 1, $lastRange.End - 1)
$splitPoint.InsertParagraphAfter()

# The freshly created (empty) paragraph now follows.
$newPara = $d.Paragraphs.Item($lastParaIndex + 1)
$newRange = $newPara.Range

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$part1 = "4. A. Each solution somewhat meets the goals the only that works for both questions is that you need "
$part2 = "to pick 20 socks to guarantee 1 or 3 pairs."

$xml = '<w:p ' + $wNs + '>' +
         '<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">' + $part1 + '</w:t></w:r>' +
         '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
         '<w:r><w:t>' + $part2 + '</w:t></w:r>' +
       '</w:p>'

$newRange.InsertXML($xml)
